$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.681.75"
$ws.Range("D2").Style = $__style
$ws.Range("E2").Value = "  +3.35%  "
$__style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.604.51"
$ws.Range("D3").Style = $__style
$ws.Range("E3").Value = "  +6.05%  "
$ws.Range("E4").Value = "  +0.01%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.26"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +4.38%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "654.70"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +6.04%  "
$ws.Range("E7").Value = "  +7.45%  "
$ws.Range("E8").Value = "  +5.32%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +5.59%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.602.94"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  +6.12%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.33"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("E13").Value = "  +2.01%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.35"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +2.40%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.288.43"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  +5.91%  "
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.542.32"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("E17").Value = "  +5.15%  "
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.597.54"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = "  +5.80%  "
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  -1.25%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.53"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  +9.72%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.15"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("E22").Value = "  +6.01%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.490"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +12.47%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "511.79"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  +3.46%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000197"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  +7.94%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.69"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  +3.20%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "97.25"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  +6.37%  "
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.84"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  +7.79%  "
$ws.Range("E29").Value = "  +18.55%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.36"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  +4.27%  "
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.996"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  +3.95%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.91"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +8.85%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.565"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +5.27%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.35"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  +12.53%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "571.29"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  +3.77%  "
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  +9.55%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.930"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("E43").Value = "  +1.57%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.75"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +5.63%  "
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.77"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +0.65%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "34.13"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +39.46%  "
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.26"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  +8.37%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0419"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +3.78%  "
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.34"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +2.67%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.18"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("E51").Value = "  -7.43%  "
